# CodeSystem-BRFabricantePNI.xlsx -- "Add CBO e faz ajustes nas pages de MI e MC"
#
# Applies to the "Concepts" sheet of the BRFabricantePNI manufacturer
# CodeSystem workbook:
#   1. Corrects / updates five existing manufacturer display & definition
#      values (rows 7, 27, 48 and 51).
#   2. Appends five brand-new manufacturer concepts as rows 88-92
#      (Level/Code/Display/Definition), extending the used range from
#      A1:D87 to A1:D92.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# ---------------------------------------------------------------------
# 1. Fix existing rows
# ---------------------------------------------------------------------

# Row 7  - BUTANTAN: "FUNDACAO BUTANTAN" -> "INSTITUTO BUTANTAN"
$ws.Range("D7").Value = "INSTITUTO BUTANTAN"

# Row 27 - WYETH: short + full name both updated
$ws.Range("C27").Value = "WYETH-LTDA"
$ws.Range("D27").Value = "WYETH INDUSTRIA FARMACEUTICA LTDA"

# Row 48 - KAMADA: full name updated
$ws.Range("D48").Value = "KAMADA LTD."

# Row 51 - PFIZER-BELGICA: fix typo BELGIVA -> BELGICA
$ws.Range("D51").Value = "PFIZER MANUFACTURING BELGIUM NV - BELGICA"

# ---------------------------------------------------------------------
# 2. Append new concept rows 88-92
# ---------------------------------------------------------------------

$newRows = @(
    @{ Row = 88; Code = "44618"; Display = "BIONTECH";       Definition = "BIONTECH MANUFACTURING GMBH" },
    @{ Row = 89; Code = "44779"; Display = "MIBE";           Definition = "MIBE GMBH ARZNEIMITTEL BRECHNA" },
    @{ Row = 90; Code = "44805"; Display = "BEIJING";        Definition = "BEIJING INSTITUTE OF BIOLOGICAL PRODUCTS CO., LTD." },
    @{ Row = 91; Code = "44781"; Display = "JUBILANT";       Definition = "JUBILANT HOLLISTERSTIER LLC" },
    @{ Row = 92; Code = "45086"; Display = "PFIZER-IRLANDA"; Definition = "PFIZER IRELAND PHARMACEUTICALS" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Column A (Level) and B (Code) are numeric-looking text ("1", "44618", ...).
    # Writing them as a literal formula that evaluates to the text, then
    # collapsing the formula down to its value via copy/paste-values, stores
    # them as genuine shared-string text (t="s") instead of letting Excel's
    # value parser coerce them to numbers.
    $ws.Range("A$r").Formula = '="1"'
    $ws.Range("B$r").Formula = '="' + $entry.Code + '"'
    $ws.Range("A$r" + ":B$r").Copy()
    $ws.Range("A$r" + ":B$r").PasteSpecial(-4163)

    $ws.Range("C$r").Value = $entry.Display
    $ws.Range("D$r").Value = $entry.Definition

    # Match the existing body-row formatting (border/alignment/wrap) by
    # copying the format down from the row directly above.
    $prev = $r - 1
    $ws.Range("A$prev" + ":D$prev").Copy()
    $ws.Range("A$r" + ":D$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
